{"js": "// Replace the 100 arithmetic expressions in the single 20x5 table with\n// their updated values (order = document/table order: row-major, top-left\n// to bottom-right). Cell text is addressed by position (not by searching\n// for the old expression text) because a couple of the old expressions\n// are not unique within the table (e.g. \"98-11=\" appears twice, with two\n// different replacements), so text-search-and-replace would be ambiguous.\n\nconst newValues = [\"45-9=\", \"82+2=\", \"65-7=\", \"16+33=\", \"74-58=\", \"42-26=\", \"80-5=\", \"49+48=\", \"19-15=\", \"91-63=\", \"15+79=\", \"80-32=\", \"27+64=\", \"50+35=\", \"2+22=\", \"66+28=\", \"58-16=\", \"91-28=\", \"59+3=\", \"33-20=\", \"23+54=\", \"40+5=\", \"79-10=\", \"96-68=\", \"88+2=\", \"17+66=\", \"60-23=\", \"32+21=\", \"74-2=\", \"44+41=\", \"24+10=\", \"70-19=\", \"1+74=\", \"75+11=\", \"98-20=\", \"96-17=\", \"63-25=\", \"16-2=\", \"90-47=\", \"47+5=\", \"95-82=\", \"1+86=\", \"53-4=\", \"2+10=\", \"13+19=\", \"14+41=\", \"85-83=\", \"84+4=\", \"35+44=\", \"35+47=\", \"56-42=\", \"33+36=\", \"84-68=\", \"24+29=\", \"21+58=\", \"15-13=\", \"2+20=\", \"31+6=\", \"21+42=\", \"12+49=\", \"46+3=\", \"59-23=\", \"80-18=\", \"29-22=\", \"89-26=\", \"27+14=\", \"34-7=\", \"24+32=\", \"67+6=\", \"79+8=\", \"91-18=\", \"59-21=\", \"91-68=\", \"39-17=\", \"78-16=\", \"33+41=\", \"32+36=\", \"86-20=\", \"81-14=\", \"93-42=\", \"67-28=\", \"12+49=\", \"41+33=\", \"30+52=\", \"5+47=\", \"68+25=\", \"73+7=\", \"51+30=\", \"13+72=\", \"50-46=\", \"33+21=\", \"77-69=\", \"58+34=\", \"49-44=\", \"14+55=\", \"64+7=\", \"42+56=\", \"15+61=\", \"83+10=\", \"63-9=\"];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columns = 5;\nconst rows = table.rowCount; // 20\n\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < columns; c++) {\n    const idx = r * columns + c;\n    if (idx >= newValues.length) continue;\n    table.getCell(r, c).value = newValues[idx];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic expressions in the single 20x5 table with\n# their updated values (order = document/table order: row-major, top-left\n# to bottom-right). Cells are addressed by (row, column) position -- not by\n# Find/Replace on the old expression text -- because a couple of the old\n# expressions are not unique within the table (e.g. \"98-11=\" appears twice,\n# with two different replacements), so text search-and-replace would be\n# ambiguous.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$newValues = @(\"45-9=\", \"82+2=\", \"65-7=\", \"16+33=\", \"74-58=\", \"42-26=\", \"80-5=\", \"49+48=\", \"19-15=\", \"91-63=\", \"15+79=\", \"80-32=\", \"27+64=\", \"50+35=\", \"2+22=\", \"66+28=\", \"58-16=\", \"91-28=\", \"59+3=\", \"33-20=\", \"23+54=\", \"40+5=\", \"79-10=\", \"96-68=\", \"88+2=\", \"17+66=\", \"60-23=\", \"32+21=\", \"74-2=\", \"44+41=\", \"24+10=\", \"70-19=\", \"1+74=\", \"75+11=\", \"98-20=\", \"96-17=\", \"63-25=\", \"16-2=\", \"90-47=\", \"47+5=\", \"95-82=\", \"1+86=\", \"53-4=\", \"2+10=\", \"13+19=\", \"14+41=\", \"85-83=\", \"84+4=\", \"35+44=\", \"35+47=\", \"56-42=\", \"33+36=\", \"84-68=\", \"24+29=\", \"21+58=\", \"15-13=\", \"2+20=\", \"31+6=\", \"21+42=\", \"12+49=\", \"46+3=\", \"59-23=\", \"80-18=\", \"29-22=\", \"89-26=\", \"27+14=\", \"34-7=\", \"24+32=\", \"67+6=\", \"79+8=\", \"91-18=\", \"59-21=\", \"91-68=\", \"39-17=\", \"78-16=\", \"33+41=\", \"32+36=\", \"86-20=\", \"81-14=\", \"93-42=\", \"67-28=\", \"12+49=\", \"41+33=\", \"30+52=\", \"5+47=\", \"68+25=\", \"73+7=\", \"51+30=\", \"13+72=\", \"50-46=\", \"33+21=\", \"77-69=\", \"58+34=\", \"49-44=\", \"14+55=\", \"64+7=\", \"42+56=\", \"15+61=\", \"83+10=\", \"63-9=\")\n\n$columns = 5\n$rows = $t.Rows.Count\n\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $columns; $c++) {\n        $idx = (($r - 1) * $columns) + ($c - 1)\n        $t.Cell($r, $c).Range.Text = $newValues[$idx]\n    }\n}\n"}
